# Updated cryptos list values (Price / Volume(1h)) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.218.58"
$ws.Cells.Item(2, 5).Value = "  +0.81%  "
$ws.Cells.Item(3, 4).Value = "1.853.56"
$ws.Cells.Item(3, 5).Value = "  +1.38%  "
$ws.Cells.Item(4, 5).Value = "  -0.47%  "
$ws.Cells.Item(5, 4).Value = "'313.78"
$ws.Cells.Item(5, 5).Value = "  +0.77%  "
$ws.Cells.Item(6, 5).Value = "  -0.39%  "
$ws.Cells.Item(7, 5).Value = "  +0.51%  "
$ws.Cells.Item(8, 4).Value = "'0.3711"
$ws.Cells.Item(8, 5).Value = "  +0.25%  "
$ws.Cells.Item(9, 4).Value = "'0.07295"
$ws.Cells.Item(9, 5).Value = "  -0.63%  "
$ws.Cells.Item(10, 4).Value = "'0.8905"
$ws.Cells.Item(10, 5).Value = "  +1.64%  "
$ws.Cells.Item(11, 4).Value = "'20.09"
$ws.Cells.Item(11, 5).Value = "  +1.81%  "
$ws.Cells.Item(12, 4).Value = "'0.07866"
$ws.Cells.Item(12, 5).Value = "  -0.17%  "
$ws.Cells.Item(13, 4).Value = "1.805.36"
$ws.Cells.Item(13, 5).Value = "  -0.86%  "
$ws.Cells.Item(14, 4).Value = "'5.403"
$ws.Cells.Item(14, 5).Value = "  +1.26%  "
$ws.Cells.Item(15, 4).Value = "'6.519"
$ws.Cells.Item(15, 5).Value = "  -0.43%  "
$ws.Cells.Item(16, 4).Value = "'91.26"
$ws.Cells.Item(16, 5).Value = "  -0.22%  "
$ws.Cells.Item(17, 5).Value = "  -0.46%  "
$ws.Cells.Item(18, 4).Value = "'0.000008916"
$ws.Cells.Item(18, 5).Value = "  +0.86%  "
$ws.Cells.Item(19, 5).Value = "  -0.31%  "
$ws.Cells.Item(21, 4).Value = "27.248.87"
$ws.Cells.Item(21, 5).Value = "  +0.84%  "
$ws.Cells.Item(22, 4).Value = "'5.089"
$ws.Cells.Item(22, 5).Value = "  -0.26%  "
$ws.Cells.Item(23, 5).Value = "  +0.16%  "
$ws.Cells.Item(24, 4).Value = "2.064.93"
$ws.Cells.Item(24, 5).Value = "  +1.39%  "
$ws.Cells.Item(25, 4).Value = "'1.954"
$ws.Cells.Item(25, 5).Value = "  +5.51%  "
$ws.Cells.Item(26, 4).Value = "'151.43"
$ws.Cells.Item(26, 5).Value = "  -0.53%  "
$ws.Cells.Item(27, 5).Value = "  -0.27%  "
$ws.Cells.Item(28, 5).Value = "  +0.14%  "
$ws.Cells.Item(29, 4).Value = "'115.88"
$ws.Cells.Item(29, 5).Value = "  +0.23%  "
$ws.Cells.Item(30, 4).Value = "'5.042"
$ws.Cells.Item(30, 5).Value = "  -1.32%  "
$ws.Cells.Item(31, 5).Value = "  -0.56%  "
$ws.Cells.Item(32, 4).Value = "'3.140"
$ws.Cells.Item(32, 5).Value = "  +6.00%  "
$ws.Cells.Item(33, 4).Value = "'0.7722"
$ws.Cells.Item(33, 5).Value = "  +5.98%  "
$ws.Cells.Item(35, 4).Value = "'4.516"
$ws.Cells.Item(35, 5).Value = "  +1.78%  "
$ws.Cells.Item(36, 4).Value = "'2.712"
$ws.Cells.Item(36, 5).Value = "  +10.41%  "
$ws.Cells.Item(37, 4).Value = "'1.111"
$ws.Cells.Item(37, 5).Value = "  +3.15%  "
$ws.Cells.Item(38, 5).Value = "  -0.13%  "
$ws.Cells.Item(39, 5).Value = "  +0.03%  "
$ws.Cells.Item(40, 4).Value = "'2.945"
$ws.Cells.Item(40, 5).Value = "  -0.33%  "
$ws.Cells.Item(41, 4).Value = "'7.054"
$ws.Cells.Item(41, 5).Value = "  -0.60%  "
$ws.Cells.Item(42, 4).Value = "'0.5126"
$ws.Cells.Item(42, 5).Value = "  -0.59%  "
$ws.Cells.Item(43, 4).Value = "'0.1625"
$ws.Cells.Item(43, 5).Value = "  -0.01%  "
$ws.Cells.Item(44, 4).Value = "'8.487"
$ws.Cells.Item(44, 5).Value = "  +4.07%  "
$ws.Cells.Item(45, 4).Value = "'0.4789"
$ws.Cells.Item(45, 5).Value = "  -0.69%  "
$ws.Cells.Item(46, 4).Value = "'10.40"
$ws.Cells.Item(46, 5).Value = "  +1.81%  "
$ws.Cells.Item(47, 5).Value = "  -0.45%  "
$ws.Cells.Item(48, 4).Value = "'102.66"
$ws.Cells.Item(48, 5).Value = "  +0.72%  "
$ws.Cells.Item(49, 5).Value = "  +1.24%  "
$ws.Cells.Item(50, 4).Value = "'0.06199"
$ws.Cells.Item(50, 5).Value = "  -0.09%  "
$ws.Cells.Item(51, 4).Value = "'65.33"
$ws.Cells.Item(51, 5).Value = "  +0.80%  "
